$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue "D2" "27.353.41"
$ws.Range("E2").Value = "  -2.97%  "
Set-TextValue "D3" "1.857.60"
$ws.Range("E3").Value = "  -3.03%  "
Set-TextValue "D4" "1.001"
$ws.Range("E4").Value = "  +0.15%  "
Set-TextValue "D5" "328.07"
$ws.Range("E5").Value = "  +0.21%  "
Set-TextValue "D6" "1.001"
$ws.Range("E6").Value = "  +0.19%  "
Set-TextValue "D7" "0.4550"
$ws.Range("E7").Value = "  -2.76%  "
Set-TextValue "D8" "0.3917"
$ws.Range("E8").Value = "  -2.28%  "
Set-TextValue "D9" "47.52"
$ws.Range("E9").Value = "  -10.64%  "
Set-TextValue "D10" "0.07925"
$ws.Range("E10").Value = "  -5.65%  "
$ws.Range("E11").Value = "  -3.07%  "
Set-TextValue "D12" "21.46"
$ws.Range("E12").Value = "  -2.96%  "
Set-TextValue "D13" "1.860.46"
$ws.Range("E13").Value = "  -2.49%  "
Set-TextValue "D14" "5.916"
$ws.Range("E14").Value = "  -2.47%  "
Set-TextValue "D15" "7.153"
$ws.Range("E15").Value = "  -3.70%  "
$ws.Range("E16").Value = "  +0.00%  "
$ws.Range("B17").Value = "Litecoin"
$ws.Range("C17").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue "D17" "86.40"
$ws.Range("E17").Value = "  -3.58%  "
$ws.Range("B18").Value = "TRON"
$ws.Range("C18").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextValue "D18" "0.06616"
$ws.Range("E18").Value = "  +0.47%  "
Set-TextValue "D19" "0.00001028"
$ws.Range("E19").Value = "  -3.00%  "
Set-TextValue "D20" "17.20"
$ws.Range("E20").Value = "  -4.15%  "
$ws.Range("E21").Value = "  +0.20%  "
Set-TextValue "D22" "5.493"
$ws.Range("E22").Value = "  -3.88%  "
Set-TextValue "D23" "27.355.03"
$ws.Range("E23").Value = "  -2.93%  "
Set-TextValue "D24" "10.92"
$ws.Range("E24").Value = "  -3.68%  "
Set-TextValue "D25" "2.295"
$ws.Range("E25").Value = "  +0.86%  "
Set-TextValue "D26" "2.080.11"
$ws.Range("E26").Value = "  -2.52%  "
Set-TextValue "D27" "154.62"
$ws.Range("E27").Value = "  +0.84%  "
Set-TextValue "D28" "20.00"
$ws.Range("E28").Value = "  +0.00%  "
Set-TextValue "D29" "2.072"
$ws.Range("E29").Value = "  -2.69%  "
Set-TextValue "D30" "5.463"
$ws.Range("E30").Value = "  -4.35%  "
Set-TextValue "D31" "121.27"
$ws.Range("E31").Value = "  -1.61%  "
Set-TextValue "D32" "0.9519"
$ws.Range("E32").Value = "  -2.39%  "
Set-TextValue "D33" "0.09376"
$ws.Range("E33").Value = "  -2.24%  "
Set-TextValue "D34" "1.456"
$ws.Range("E34").Value = "  +0.72%  "
Set-TextValue "D35" "3.588"
$ws.Range("E35").Value = "  -1.07%  "
Set-TextValue "D36" "5.265"
$ws.Range("E36").Value = "  -5.01%  "
Set-TextValue "D37" "0.06037"
$ws.Range("E37").Value = "  -1.89%  "
Set-TextValue "D38" "0.02229"
$ws.Range("E38").Value = "  -3.00%  "
Set-TextValue "D39" "1.219"
$ws.Range("E39").Value = "  -2.26%  "
Set-TextValue "D40" "8.054"
$ws.Range("E40").Value = "  -8.82%  "
$ws.Range("E41").Value = "  +0.14%  "
Set-TextValue "D42" "0.5922"
$ws.Range("E42").Value = "  -3.41%  "
Set-TextValue "D43" "0.1884"
$ws.Range("E43").Value = "  -0.95%  "
Set-TextValue "D44" "10.16"
$ws.Range("E44").Value = "  -7.77%  "
$ws.Range("E45").Value = "  -1.30%  "
$ws.Range("E46").Value = "  -3.97%  "
Set-TextValue "D47" "12.07"
$ws.Range("E47").Value = "  -5.39%  "
Set-TextValue "D48" "3.393"
$ws.Range("E48").Value = "  -1.14%  "
Set-TextValue "D49" "1.918"
$ws.Range("E49").Value = "  -5.31%  "
Set-TextValue "D50" "0.06746"
$ws.Range("E50").Value = "  -1.69%  "
Set-TextValue "D51" "108.26"
$ws.Range("E51").Value = "  -1.29%  "
